# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_device")

# Update the "device id" values in column A for the specified rows:
# Rows that had 10002 -> 10003
$rowsFrom10002 = @(3, 23, 43, 63, 83)
foreach ($r in $rowsFrom10002) {
    $ws.Cells.Item($r, 1).Value = 10003
}

# Rows that had 10005 -> 10003
$rowsFrom10005 = @(105, 114, 123, 132, 141)
foreach ($r in $rowsFrom10005) {
    $ws.Cells.Item($r, 1).Value = 10003
}

# Update the view/selection state: scroll back to top-left and select from A162
# down to the end of the sheet (whole rows below the data), matching the
# final sheetView/selection recorded in the saved file.
$ws.Activate()
$ws.Range("A162:XFD1048576").Select()
